$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.720.65'
$ws.Range("E2").Value = '  +2.41%  '
$ws.Range("D3").Value = '2.160.16'
$ws.Range("E3").Value = '  +2.80%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.11'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("E6").Value = '  +1.75%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.26'
$ws.Range("E7").Value = '  +1.70%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.391'
$ws.Range("E9").Value = '  +0.67%  '
$ws.Range("E10").Value = '  -0.08%  '
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.90'
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").Value = '2.479.29'
$ws.Range("E13").Value = '  +2.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.82'
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.49'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").Value = '2.146.18'
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("D18").Value = '39.623.99'
$ws.Range("E18").Value = '  +2.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.63'
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("E20").Value = '  -1.01%  '
$ws.Range("D21").Value = '0.0₃0846'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.41'
$ws.Range("E22").Value = '  +0.93%  '
$ws.Range("E24").Value = '  +2.35%  '
$ws.Range("E25").Value = '  -8.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '172.44'
$ws.Range("E26").Value = '  +1.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.55'
$ws.Range("E27").Value = '  -1.02%  '
$ws.Range("E28").Value = '  +2.27%  '
$ws.Range("E29").Value = '  +2.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.79'
$ws.Range("E31").Value = '  +5.90%  '
$ws.Range("E32").Value = '  +1.46%  '
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("E34").Value = '  -1.52%  '
$ws.Range("E35").Value = '  -3.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0616'
$ws.Range("E36").Value = '  +0.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.67'
$ws.Range("E37").Value = '  +4.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.41'
$ws.Range("E38").Value = '  +1.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.10'
$ws.Range("E39").Value = '  +22.62%  '
$ws.Range("E40").Value = '  +0.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '102.66'
$ws.Range("E41").Value = '  +1.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.78'
$ws.Range("E42").Value = '  -1.03%  '
$ws.Range("E43").Value = '  -0.50%  '
$ws.Range("D44").Value = '1.515.07'
$ws.Range("E44").Value = '  -0.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.21'
$ws.Range("E45").Value = '  +1.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.86'
$ws.Range("E46").Value = '  +1.02%  '
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0918'
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("E49").Value = '  +0.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '49.99'
$ws.Range("E50").Value = '  +8.62%  '
$ws.Range("E51").Value = '  +1.14%  '
